# Applies the "tervezet.docx" edit:
#   - Bejelentkezésnél ... paragraph -> font color FFC000 (gold)
#   - Téma hozzáadásánál ... paragraph -> font color FFC000 (gold)
#   - Főmenü gomb paragraph -> font color FFFF00 (yellow)
#   - "Ügyintézés alatt és lezárt ügyek törlése" paragraph -> removed entirely
#   - "Sötét mód hozzáadása" -> text changed to "Sötét mód " + "fixálása"
#     (as two separate runs, second word swapped)

$d = $word.ActiveDocument

function Find-ParagraphByPrefix($doc, $prefix) {
    foreach ($para in $doc.Paragraphs) {
        if ($para.Range.Text.StartsWith($prefix)) {
            return $para
        }
    }
    return $null
}

# wdColor values are stored as 0x00BBGGRR, so RGB hex FFC000 / FFFF00 map to:
$goldColor = 49407   # FFC000
$yellowColor = 65535 # FFFF00

# 1) "Bejelentkezésnél először kelljen..." -> gold
$p1 = Find-ParagraphByPrefix $d "Bejelentkezésnél"
$p1.Range.Font.Color = $goldColor

# 2) "Téma hozzáadásánál is kelljen..." -> gold
$p2 = Find-ParagraphByPrefix $d "Téma hozzáadásánál"
$p2.Range.Font.Color = $goldColor

# 3) "Főmenü gomb" -> yellow
$p3 = Find-ParagraphByPrefix $d "Főmenü gomb"
$p3.Range.Font.Color = $yellowColor

# 4) Remove the "Ügyintézés alatt és lezárt ügyek törlése" bullet entirely
$p4 = Find-ParagraphByPrefix $d "Ügyintézés alatt és lezárt ügyek törlése"
$p4.Range.Delete()

# 5) "Sötét mód hozzáadása" -> "Sötét mód " + "fixálása" (two runs)
$p5 = Find-ParagraphByPrefix $d "Sötét mód hozzáadása"
$fullRange = $p5.Range
$textRange = $d.Range($fullRange.Start, $fullRange.End - 1)
$textRange.InsertXML('<w:p><w:pPr><w:pStyle w:val="Listaszerbekezds"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Sötét mód </w:t></w:r><w:r><w:t>fixálása</w:t></w:r></w:p>')
